# Delete the "LHE" (Lahore, Pakistan) row from the colo data sheet.
# This shifts every following row up by one, which matches the target diff
# (e.g. row 219 "PKX" becomes row 218, ..., row 332 "YHZ" becomes row 331),
# and reduces the sheet's used range from A1:H332 to A1:H331.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(218).Delete()
